$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Populate row 27 with the new activity log entry (use raw Excel serial
# numbers so the existing cell number-format styles are preserved)
$ws.Cells.Item(27, 2).Value = 6977
$ws.Cells.Item(27, 3).Value = 43924
$ws.Cells.Item(27, 4).Value = 0.92708333333333337
$ws.Cells.Item(27, 5).Value = 0.96875
$ws.Cells.Item(27, 7).Value = "Fixed ArithUnit.vhd and Adder.vhd such that both .do scripts work (previosuly they do not work as pointed out by my teammates). Each script works if and only if all ports are defined properly"

# Update the selected cell to match the new active cell
$ws.Range("B30").Select()
